# Expand the academic-editor "Response:" paragraph (criteria re: utility,
# validation, availability) per the reviewer-response revision.
#
# The target paragraph is the last paragraph before the "Reviewer 1"
# heading: a short "Response:" paragraph that follows the "To improve the
# chances of manuscript suitability for PLOS ONE..." bullet. We locate it
# by scanning paragraphs for that exact text, immediately followed by a
# paragraph whose text is "Reviewer 1".

$d = $word.ActiveDocument

$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Response:") {
        if ($i -lt $count) {
            $nextTxt = $d.Paragraphs($i + 1).Range.Text.TrimEnd([char]13, [char]7)
            if ($nextTxt -eq "Reviewer 1") {
                $targetIndex = $i
                break
            }
        }
    }
}

if ($targetIndex -eq -1) {
    Write-Host "ERROR: could not locate target Response paragraph"
} else {
    $p = $d.Paragraphs($targetIndex)

    # --- Extend the plain ":" run into the full response sentence -------
    $colonRange = $p.Range
    $found = $colonRange.Find.Execute(":")
    if ($found) {
        $colonRange.Text = ": We have verified that our manuscript adheres to the utility, validation, and availability criteria for the article type. Below are descriptions of how these criteria are met:"
    } else {
        Write-Host "ERROR: could not find trailing colon run"
    }

    # --- Insert the Utility / Validation / Availability paragraphs ------
    $p.Range.InsertParagraphAfter()
    $pUtility = $d.Paragraphs($targetIndex + 1)
    $pUtility.Range.InsertAfter("Utility: The MassWateR package provides needed utility for the community, as established in the introductory text. both in the original submission and the additions provided in the revision. In particular, the package provides a repeatable and efficient means for generating QC reports to ensure the water quality data are of sufficient accuracy and precision needed for integration into larger databases (i.e., WQX). As noted, we are unaware of any existing tools that provide this functionality.")
    $pUtility.Style = "BodyText"

    $pUtility.Range.InsertParagraphAfter()
    $pValidation = $d.Paragraphs($targetIndex + 2)
    $pValidation.Range.InsertAfter("Validation: MassWateR achieves its intended uses as verified by adoption by the larger community of practice, detailed in the text. In addition, the package has been downloaded over 4000 times since its availability on CRAN in January 2023. This provides a robust indication that those outside of our community of practice are also applying MassWateR for its intended use. Finally, we describe the possible applications for the package in detail in the manuscript (e.g., QC report generation including a figure showing a portion of the report, submission to WQX, etc.).")
    $pValidation.Style = "BodyText"

    $pValidation.Range.InsertParagraphAfter()
    $pAvailability = $d.Paragraphs($targetIndex + 3)
    $pAvailability.Range.InsertAfter("Availability: The software is entirely open -source and available for download through ~CRAN~, ~GITHUB~, and ~RUNIVERSE~ Numerous links and citations are included throughout the text. The license applied to MassWateR is CC0 1.0, making the code available to the world-wide public domain.")
    $pAvailability.Style = "BodyText"

    # --- Italicize the three bold labels ---------------------------------
    $r = $pUtility.Range
    if ($r.Find.Execute("Utility:")) { $r.Italic = 1 }

    $r = $pValidation.Range
    if ($r.Find.Execute("Validation")) { $r.Italic = 1 }

    $r = $pAvailability.Range
    if ($r.Find.Execute("Availability")) { $r.Italic = 1 }

    # --- Turn the placeholder tokens into hyperlinks ----------------------
    $r = $pAvailability.Range
    if ($r.Find.Execute("~CRAN~")) {
        $d.Hyperlinks.Add($r, "https://cran.r-project.org/package=MassWateR", $null, $null, "CRAN") | Out-Null
    }

    $r = $pAvailability.Range
    if ($r.Find.Execute("~GITHUB~")) {
        $d.Hyperlinks.Add($r, "https://github.com/massbays-tech/MassWateR", $null, $null, "GitHub") | Out-Null
    }

    $r = $pAvailability.Range
    if ($r.Find.Execute("~RUNIVERSE~")) {
        $d.Hyperlinks.Add($r, "https://massbays-tech.r-universe.dev/MassWateR", $null, $null, "R-Universe") | Out-Null
    }

    Write-Host "Updated paragraph: $($p.Range.Text)"
    Write-Host "Inserted paragraph: $($pUtility.Range.Text)"
    Write-Host "Inserted paragraph: $($pValidation.Range.Text)"
    Write-Host "Inserted paragraph: $($pAvailability.Range.Text)"
}
